# WRESBAL FRED data refresh: append the two newest weekly observations
# to the "Data" sheet and refresh the FRED series metadata on
# "SeriesInfo" to match the latest pull.

$wb = $excel.ActiveWorkbook

# --- Sheet "Data": append two new weekly observations ---
$data = $wb.Worksheets.Item("Data")

# Row 110 - copy the date-column formatting from the last existing row
# (style includes the YYYY-MM-DD number format + bold/centered/bordered
# header look used throughout column A) before writing the new values.
$data.Cells.Item(109, 1).Copy($data.Cells.Item(110, 1))
$data.Cells.Item(110, 1).Value = 45231
$data.Cells.Item(110, 2).Value = 3267.216

# Row 111
$data.Cells.Item(109, 1).Copy($data.Cells.Item(111, 1))
$data.Cells.Item(111, 1).Value = 45238
$data.Cells.Item(111, 2).Value = 3328.908

# --- Sheet "SeriesInfo": refresh metadata fields from the new FRED pull ---
$info = $wb.Worksheets.Item("SeriesInfo")

# These hold plain date-looking text (e.g. "2023-11-15"); force Text
# formatting first so they are stored as strings rather than being
# auto-converted to date serial numbers, then drop the format override
# so the cell stays unstyled like its neighbours.
$info.Range("B3").NumberFormat = "@"
$info.Range("B3").Value = "2023-11-15"
$info.Range("B3").ClearFormats()

$info.Range("B4").NumberFormat = "@"
$info.Range("B4").Value = "2023-11-15"
$info.Range("B4").ClearFormats()

$info.Range("B7").NumberFormat = "@"
$info.Range("B7").Value = "2023-11-08"
$info.Range("B7").ClearFormats()

# Timestamp string with a UTC-offset suffix isn't auto-parsed as a date,
# so it can be set directly.
$info.Range("B14").Value = "2023-11-09 15:39:01-06"

# Popularity score - stays numeric.
$info.Range("B15").Value = 73
